# Updates cryptos list data (Price and Volume(1h) columns) per upstream scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.269.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.590.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.90%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "509.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.602.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.86%  "
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.047.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.283.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.598.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.66%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.421"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0840"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.859"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.13%  "
$ws.Range("E38").Value = "  -4.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.840"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "295.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.82%  "
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.618"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0553"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.11%  "
$ws.Range("E49").Value = "  -2.17%  "
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.992.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.72%  "
